$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Anchors (unchanged cells) used only to clone exact cell formatting via
# PasteSpecial so that type-changing cells land on the same style index
# the target workbook uses, instead of a freshly synthesized duplicate.
#   D14 -> style 13 (General/text, "N/A"-style cell)
#   N14 -> style 14 (#,##0.0;"-"#,##0.0  percent-change cell)
#   C15 -> style 15 (#,##0  integer-count cell)
$xlPasteFormats = -4122

# --- Header text updates (rich-text cells, set whole value) ---
$ws.Range("A8").Value = "Volume 32   Number  2"
$ws.Range("C9").Value = "Report Covering the Week  1/6/2025  Through  1/12/2025"

# --- Row 15 (Rape) ---
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 2
$ws.Range("L15").Value = 100
$ws.Range("N14").Copy()
$ws.Range("L15").PasteSpecial($xlPasteFormats)
$ws.Range("M15").Value = 100

# --- Row 16 (Robbery) ---
$ws.Range("C16").Value = 3
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 11
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 6
$ws.Range("J16").Value = 5
$ws.Range("K16").Value = 20
$ws.Range("L16").Value = -14.285714285714
$ws.Range("M16").Value = -53.846153846153
$ws.Range("N16").Value = -83.783783783783

# --- Row 17 (Fel. Assault) ---
$ws.Range("C17").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("C17").PasteSpecial($xlPasteFormats)
$ws.Range("E17").Value = -100
$ws.Range("G17").Value = 27
$ws.Range("H17").Value = -70.370370370370
$ws.Range("J17").Value = 11
$ws.Range("K17").Value = -63.636363636363
$ws.Range("L17").Value = -42.857142857142
$ws.Range("M17").Value = -42.857142857142
$ws.Range("N17").Value = -84

# --- Row 18 (Burglary) ---
$ws.Range("C18").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("C18").PasteSpecial($xlPasteFormats)
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 4
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = -73.333333333333
$ws.Range("J18").Value = 7
$ws.Range("K18").Value = -85.714285714285
$ws.Range("L18").Value = -88.888888888888
$ws.Range("M18").Value = -85.714285714285
$ws.Range("N18").Value = -93.75

# --- Row 19 (Gr. Larceny) ---
$ws.Range("C19").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("C19").PasteSpecial($xlPasteFormats)
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = -100
$ws.Range("F19").Value = 11
$ws.Range("G19").Value = 16
$ws.Range("H19").Value = -31.25
$ws.Range("J19").Value = 11
$ws.Range("K19").Value = -63.636363636363
$ws.Range("L19").Value = -80
$ws.Range("M19").Value = -33.333333333333
$ws.Range("N19").Value = -71.428571428571

# --- Row 20 (G.L.A.) ---
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 50
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 33.333333333333
$ws.Range("I20").Value = 4
$ws.Range("J20").Value = 3
$ws.Range("K20").Value = 33.333333333333
$ws.Range("L20").Value = 33.333333333333
$ws.Range("M20").Value = 100
$ws.Range("N20").Value = -84.615384615384

# --- Row 21 (TOTAL) ---
$ws.Range("C21").Value = 7
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = -70.833333333333
$ws.Range("F21").Value = 44
$ws.Range("G21").Value = 76
$ws.Range("H21").Value = -42.105263157894
$ws.Range("I21").Value = 21
$ws.Range("J21").Value = 37
$ws.Range("K21").Value = -43.243243243243
$ws.Range("L21").Value = -55.319148936170
$ws.Range("M21").Value = -43.243243243243
$ws.Range("N21").Value = -82.644628099173

# --- Row 22 (Transit) ---
$ws.Range("C22").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("C22").PasteSpecial($xlPasteFormats)
$ws.Range("D22").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("D22").PasteSpecial($xlPasteFormats)
$ws.Range("E22").Value = "'***.*"
$ws.Range("D14").Copy()
$ws.Range("E22").PasteSpecial($xlPasteFormats)
$ws.Range("M22").Value = 0
$ws.Range("N14").Copy()
$ws.Range("M22").PasteSpecial($xlPasteFormats)

# --- Row 23 (Housing) ---
$ws.Range("C23").Value = 1
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 3
$ws.Range("J23").Value = 2
$ws.Range("K23").Value = 50
$ws.Range("L23").Value = 50
$ws.Range("M23").Value = 50

# --- Row 24 (Petit Larceny) ---
$ws.Range("C24").Value = 13
$ws.Range("D24").Value = 17
$ws.Range("E24").Value = -23.529411764705
$ws.Range("F24").Value = 52
$ws.Range("G24").Value = 51
$ws.Range("H24").Value = 1.960784313725
$ws.Range("I24").Value = 19
$ws.Range("J24").Value = 21
$ws.Range("K24").Value = -9.523809523809
$ws.Range("L24").Value = -29.629629629629
$ws.Range("M24").Value = -5

# --- Row 25 (Retail Theft) ---
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 1
$ws.Range("C15").Copy()
$ws.Range("D25").PasteSpecial($xlPasteFormats)
$ws.Range("E25").Value = 300
$ws.Range("N14").Copy()
$ws.Range("E25").PasteSpecial($xlPasteFormats)
$ws.Range("I25").Value = 6
$ws.Range("J25").Value = 1
$ws.Range("C15").Copy()
$ws.Range("J25").PasteSpecial($xlPasteFormats)
$ws.Range("K25").Value = 500
$ws.Range("N14").Copy()
$ws.Range("K25").PasteSpecial($xlPasteFormats)
$ws.Range("L25").Value = -14.285714285714

# --- Row 26 (Misd. Assault) ---
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = -57.142857142857
$ws.Range("F26").Value = 34
$ws.Range("G26").Value = 29
$ws.Range("H26").Value = 17.241379310344
$ws.Range("I26").Value = 13
$ws.Range("J26").Value = 11
$ws.Range("K26").Value = 18.181818181818
$ws.Range("L26").Value = 30
$ws.Range("M26").Value = -59.375

# --- Row 27 (UCR Rape*) ---
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 2
$ws.Range("L27").Value = 0
$ws.Range("N14").Copy()
$ws.Range("L27").PasteSpecial($xlPasteFormats)

# --- Row 28 (Other Sex Crimes) ---
$ws.Range("C28").Value = 2
$ws.Range("C15").Copy()
$ws.Range("C28").PasteSpecial($xlPasteFormats)
$ws.Range("D28").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("D28").PasteSpecial($xlPasteFormats)
$ws.Range("E28").Value = "'***.*"
$ws.Range("D14").Copy()
$ws.Range("E28").PasteSpecial($xlPasteFormats)
$ws.Range("I28").Value = 2
$ws.Range("C15").Copy()
$ws.Range("I28").PasteSpecial($xlPasteFormats)
$ws.Range("K28").Value = 0

# --- Row 29 (Shooting Vic.) ---
$ws.Range("G29").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("G29").PasteSpecial($xlPasteFormats)
$ws.Range("H29").Value = "'***.*"
$ws.Range("D14").Copy()
$ws.Range("H29").PasteSpecial($xlPasteFormats)

# --- Row 30 (Shooting Inc.) ---
$ws.Range("G30").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("G30").PasteSpecial($xlPasteFormats)
$ws.Range("H30").Value = "'***.*"
$ws.Range("D14").Copy()
$ws.Range("H30").PasteSpecial($xlPasteFormats)

# --- Row 40 (historical Rape) ---
$ws.Range("J40").Value = 12
$ws.Range("K40").Value = -61.290322580645
$ws.Range("L40").Value = -58.620689655172
$ws.Range("M40").Value = -84.810126582278
$ws.Range("N40").Value = -76.470588235294

# --- Row 46 (historical TOTAL) ---
$ws.Range("J46").Value = 875
$ws.Range("K46").Value = -47.916666666666
$ws.Range("L46").Value = -56.052235057759
$ws.Range("M46").Value = -78.239243969161
$ws.Range("N46").Value = -80.936819172113
